# Agrega soporte para rangos de NUM, STRING y SWITCH
#
# - New named ranges: PUE.NUM.legajos (A3:A5) and PUE.SWITCH.booleanos (I9:I12)
# - New "LEGAJOS" column (A2:A5) mirroring the existing Nombre/Apellido/Edad table
# - New "Maquina 1/2/3" matrix (I2:L5) with Velocidad/Uso/Anios rows
# - New boolean column (I9:I12) next to the existing Porcentajes slider
# - Border/centered formatting carried over to G11, D14 and E14
# - E14 ("Analizar?") flipped from TRUE to FALSE
# - Selection moved to N6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Defined names (PUE.NUM.* / PUE.SWITCH.* ranges) ---
$wb.Names.Add("PUE.NUM.legajos", '=Sheet1!$A$3:$A$5')
$wb.Names.Add("PUE.SWITCH.booleanos", '=Sheet1!$I$9:$I$12')

# --- Column A: new "LEGAJOS" header + legajo numbers, styled like column B ---
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A2").Value = "LEGAJOS"

$ws.Range("B3").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)
$ws.Range("A3").Value = 25407
$ws.Range("A4").Value = 26549
$ws.Range("A5").Value = 23403

# --- Row 2 (I2:L2): machine headers, styled like the bordered/centered E3 cell ---
$ws.Range("E3").Copy()
$ws.Range("I2:L2").PasteSpecial(-4122)
$ws.Range("J2").Value = "Maquina 1"
$ws.Range("K2").Value = "Maquina 2"
$ws.Range("L2").Value = "Maquina 3"

# --- Rows 3-5 (I:L): Velocidad/Uso/Anios data matrix, styled like bordered D3 ---
$ws.Range("D3").Copy()
$ws.Range("I3:L5").PasteSpecial(-4122)

$ws.Range("I3").Value = "Velocidad"
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 200
$ws.Range("L3").Value = 700

$ws.Range("I4").Value = "Uso"
$ws.Range("J4").Value = 0.9
$ws.Range("K4").Value = 0.3
$ws.Range("L4").Value = 0.5

$ws.Range("I5").Value = "Anios"
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 2

# --- Column I rows 9-12: boolean switches, styled like B9 ---
$ws.Range("B9").Copy()
$ws.Range("I9:I12").PasteSpecial(-4122)
$ws.Range("I9").Value = $false
$ws.Range("I10").Value = $false
$ws.Range("I11").Value = $true
$ws.Range("I12").Value = $false

# --- G11 picks up the same bordered style as B11 ---
$ws.Range("B11").Copy()
$ws.Range("G11").PasteSpecial(-4122)

# --- D14 / E14 also get the bordered style, and the switch flips to FALSE ---
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E14").Value = $false

# --- Move the active selection to N6 ---
$ws.Range("N6").Select()
